$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell corrections (row positions unaffected by the later row deletions) ---
$ws.Range("C6").Value = 15.1
$ws.Range("C8").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("C14").ClearContents()
$ws.Range("C17").Value = 11.2
$ws.Range("C18").Value = 11.5
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("C23").Value = 12.2

# --- Remove the "RM 232" record (row 26) entirely; rows below shift up by one ---
$ws.Rows(26).Delete()

# --- Remove the "SC 92" record, now sitting at row 27 after the previous delete ---
$ws.Rows(27).Delete()

# --- Final touch-ups on the rows that settled into their new positions ---
# Row 27 is now "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()

# Row 29 is now "SC 119"
$ws.Range("B29").ClearContents()

# Row 32 is now "SC 193"
$ws.Range("B32").ClearContents()
